# Applies updated currentAveragePrice / LevePrice / LeveProfit figures
# (columns H-N) across several sheets, per the scheduled market-data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 735.9091
$ws.Range("I12").Value = 999
$ws.Range("J12").Value = 585.5714
$ws.Range("K12").Value = 999
$ws.Range("L12").Value = 585.5714
$ws.Range("M12").Value = -829
$ws.Range("N12").Value = -925.5714

$ws.Range("H18").Value = 970.7143
$ws.Range("I18").Value = 970.7143
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 970.7143
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -686.7143

$ws.Range("H33").Value = 402
$ws.Range("I33").Value = 427.33334
$ws.Range("J33").Value = 98
$ws.Range("K33").Value = 427.33334
$ws.Range("L33").Value = 98
$ws.Range("M33").Value = -198.33334
$ws.Range("N33").Value = -556

$ws.Range("H116").Value = 0
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("M116").ClearContents()

$ws.Range("H138").Value = 5219.1816
$ws.Range("I138").Value = 4557
$ws.Range("J138").Value = 5397.4614
$ws.Range("K138").Value = 13671
$ws.Range("L138").Value = 16192.3842
$ws.Range("M138").Value = -8531
$ws.Range("N138").Value = -26472.3842

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4148.5
$ws.Range("I2").Value = 4148.5
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 4148.5
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -4035.5
$ws.Range("N2").ClearContents()

$ws.Range("H4").Value = 879
$ws.Range("I4").Value = 879
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 879
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -763

$ws.Range("H5").Value = 900
$ws.Range("I5").Value = 900
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 900
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -788

$ws.Range("H32").Value = 9086.888999999999
$ws.Range("I32").Value = 8556.125
$ws.Range("J32").Value = 13333
$ws.Range("K32").Value = 8556.125
$ws.Range("L32").Value = 13333
$ws.Range("M32").Value = -8269.125
$ws.Range("N32").Value = -13907

$ws.Range("H45").Value = 1256.4
$ws.Range("I45").Value = 1256.4
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 1256.4
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -879.4000000000001
$ws.Range("N45").ClearContents()

$ws.Range("H116").Value = 4148.5
$ws.Range("I116").Value = 4148.5
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 4148.5
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = -1854.5
$ws.Range("N116").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4148.5
$ws.Range("I3").Value = 4148.5
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 4148.5
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -4034.5
$ws.Range("N3").ClearContents()

$ws.Range("H4").Value = 900
$ws.Range("I4").Value = 900
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 900
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -785

$ws.Range("H15").Value = 48333.332
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 48333.332
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 48333.332
$ws.Range("N15").Value = -48787.332

$ws.Range("H102").Value = 27633.6
$ws.Range("I102").Value = 22139
$ws.Range("J102").Value = 49612
$ws.Range("K102").Value = 22139
$ws.Range("L102").Value = 49612
$ws.Range("M102").Value = -18894
$ws.Range("N102").Value = -56102

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 734.4
$ws.Range("I22").Value = 734.4
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 734.4
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -384.4
$ws.Range("N22").ClearContents()

$ws.Range("H41").Value = 22975.166
$ws.Range("I41").Value = 7000
$ws.Range("J41").Value = 38950.332
$ws.Range("K41").Value = 7000
$ws.Range("L41").Value = 38950.332
$ws.Range("M41").Value = -6572
$ws.Range("N41").Value = -39806.332

$ws.Range("H87").Value = 0
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()

$ws.Range("H90").Value = 0
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("K90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H108").Value = 415.8
$ws.Range("I108").Value = 415.8
$ws.Range("J108").Value = 0
$ws.Range("K108").Value = 1247.4
$ws.Range("L108").Value = 0
$ws.Range("M108").Value = 1632.6

$ws.Range("H109").Value = 475.66666
$ws.Range("I109").Value = 475.66666
$ws.Range("J109").Value = 0
$ws.Range("K109").Value = 1426.99998
$ws.Range("L109").Value = 0
$ws.Range("M109").Value = -386.9999800000001

$ws.Range("H120").Value = 9013.5
$ws.Range("I120").Value = 3027
$ws.Range("J120").Value = 15000
$ws.Range("K120").Value = 9081
$ws.Range("L120").Value = 45000
$ws.Range("M120").Value = -4243
$ws.Range("N120").Value = -54676

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 177.42105
$ws.Range("I2").Value = 228
$ws.Range("J2").Value = 90.71429000000001
$ws.Range("K2").Value = 228
$ws.Range("L2").Value = 90.71429000000001
$ws.Range("M2").Value = -115
$ws.Range("N2").Value = -316.71429

$ws.Range("H43").Value = 20109.125
$ws.Range("I43").Value = 6713.5
$ws.Range("J43").Value = 33504.75
$ws.Range("K43").Value = 6713.5
$ws.Range("L43").Value = 33504.75
$ws.Range("M43").Value = -6562.5
$ws.Range("N43").Value = -33806.75

$ws.Range("H114").Value = 0
$ws.Range("I114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("K114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()

$ws.Range("H138").Value = 0
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1000
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 1000
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 1000
$ws.Range("N22").Value = -1590

$ws.Range("H23").Value = 500
$ws.Range("I23").Value = 500
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 500
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = -270

$ws.Range("H27").Value = 1000
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 1000
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 1000
$ws.Range("N27").Value = -1214

$ws.Range("H31").Value = 4753.75
$ws.Range("I31").Value = 3206
$ws.Range("J31").Value = 7333.3335
$ws.Range("K31").Value = 3206
$ws.Range("L31").Value = 7333.3335
$ws.Range("M31").Value = -2958
$ws.Range("N31").Value = -7829.3335

$ws.Range("H46").Value = 900
$ws.Range("I46").Value = 900
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 900
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = -712

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H40").Value = 25000
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 25000
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 25000
$ws.Range("N40").Value = -25298

$ws.Range("H113").Value = 295.5
$ws.Range("I113").Value = 279.5
$ws.Range("J113").Value = 327.5
$ws.Range("K113").Value = 838.5
$ws.Range("L113").Value = 982.5
$ws.Range("M113").Value = 1331.5
$ws.Range("N113").Value = -5322.5
